$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C90) from serial date 45182 to 45184
$newValue = 45184
for ($row = 2; $row -le 90; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
